$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 4705.8066
$ws.Range("I98").Value = 6644
$ws.Range("J98").Value = 1181.8182
$ws.Range("K98").Value = 6644
$ws.Range("L98").Value = 1181.8182
$ws.Range("M98").Value = -5146
$ws.Range("N98").Value = -4177.8182
$ws.Range("H115").Value = 400
$ws.Range("I115").Value = 400
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 1200
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = 367
$ws.Range("N115").ClearContents()
$ws.Range("H116").Value = 2115.8
$ws.Range("I116").Value = 1730
$ws.Range("J116").Value = 2323.5386
$ws.Range("K116").Value = 1730
$ws.Range("L116").Value = 2323.5386
$ws.Range("M116").Value = 1712
$ws.Range("N116").Value = -9207.5386
$ws.Range("H122").Value = 4705.8066
$ws.Range("I122").Value = 6644
$ws.Range("J122").Value = 1181.8182
$ws.Range("K122").Value = 19932
$ws.Range("L122").Value = 3545.4546
$ws.Range("M122").Value = -17482
$ws.Range("N122").Value = -8445.454600000001
$ws.Range("H132").Value = 3404518.5
$ws.Range("I132").Value = 3574494.2
$ws.Range("K132").Value = 10723482.6
$ws.Range("M132").Value = -10720952.6
$ws.Range("H137").Value = 5264411
$ws.Range("I137").Value = 1476.4546
$ws.Range("K137").Value = 4429.3638
$ws.Range("M137").Value = -1879.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 694.8889
$ws.Range("I45").Value = 717.5
$ws.Range("J45").Value = 514
$ws.Range("K45").Value = 717.5
$ws.Range("L45").Value = 514
$ws.Range("M45").Value = -340.5
$ws.Range("N45").Value = -1268
$ws.Range("H61").Value = 1505.4138
$ws.Range("I61").Value = 911.5714
$ws.Range("J61").Value = 3064.25
$ws.Range("K61").Value = 911.5714
$ws.Range("L61").Value = 3064.25
$ws.Range("M61").Value = -699.5714
$ws.Range("N61").Value = -3488.25
$ws.Range("H110").Value = 2168.5293
$ws.Range("I110").Value = 2164.4546
$ws.Range("K110").Value = 2164.4546
$ws.Range("M110").Value = -119.4546
$ws.Range("H136").Value = 1505.4138
$ws.Range("I136").Value = 911.5714
$ws.Range("J136").Value = 3064.25
$ws.Range("K136").Value = 2734.7142
$ws.Range("L136").Value = 9192.75
$ws.Range("M136").Value = -184.7142000000003
$ws.Range("N136").Value = -14292.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1442.8462
$ws.Range("I107").Value = 1053.875
$ws.Range("J107").Value = 2065.2
$ws.Range("K107").Value = 1053.875
$ws.Range("L107").Value = 2065.2
$ws.Range("M107").Value = 866.125
$ws.Range("N107").Value = -5905.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1056.238
$ws.Range("I16").Value = 580.13336
$ws.Range("J16").Value = 2246.5
$ws.Range("K16").Value = 580.13336
$ws.Range("L16").Value = 2246.5
$ws.Range("M16").Value = -293.13336
$ws.Range("N16").Value = -2820.5
$ws.Range("H31").Value = 3573255.8
$ws.Range("I31").Value = 1617.122
$ws.Range("J31").Value = 13335734
$ws.Range("K31").Value = 1617.122
$ws.Range("L31").Value = 13335734
$ws.Range("M31").Value = -1322.122
$ws.Range("N31").Value = -13336324
$ws.Range("H34").Value = 3573255.8
$ws.Range("I34").Value = 1617.122
$ws.Range("J34").Value = 13335734
$ws.Range("K34").Value = 1617.122
$ws.Range("L34").Value = 13335734
$ws.Range("M34").Value = -1415.122
$ws.Range("N34").Value = -13336138
$ws.Range("H58").Value = 729.06384
$ws.Range("I58").Value = 674.2564
$ws.Range("J58").Value = 996.25
$ws.Range("K58").Value = 674.2564
$ws.Range("L58").Value = 996.25
$ws.Range("M58").Value = -471.2564
$ws.Range("N58").Value = -1402.25
$ws.Range("H107").Value = 587.5
$ws.Range("I107").Value = 584.7857
$ws.Range("J107").Value = 606.5
$ws.Range("K107").Value = 584.7857
$ws.Range("L107").Value = 606.5
$ws.Range("M107").Value = 1335.2143
$ws.Range("N107").Value = -4446.5
$ws.Range("H113").Value = 1056.238
$ws.Range("I113").Value = 580.13336
$ws.Range("J113").Value = 2246.5
$ws.Range("K113").Value = 580.13336
$ws.Range("L113").Value = 2246.5
$ws.Range("M113").Value = 1589.86664
$ws.Range("N113").Value = -6586.5
$ws.Range("H136").Value = 729.06384
$ws.Range("I136").Value = 674.2564
$ws.Range("J136").Value = 996.25
$ws.Range("K136").Value = 2022.7692
$ws.Range("L136").Value = 2988.75
$ws.Range("M136").Value = 527.2308
$ws.Range("N136").Value = -8088.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 9257.25
$ws.Range("J105").Value = 11009.667
$ws.Range("L105").Value = 33029.001
$ws.Range("N105").Value = -38271.001
$ws.Range("H110").Value = 13686.692
$ws.Range("I110").Value = 785.4
$ws.Range("J110").Value = 21750
$ws.Range("K110").Value = 2356.2
$ws.Range("L110").Value = 65250
$ws.Range("M110").Value = 1733.8
$ws.Range("N110").Value = -73430
$ws.Range("H131").Value = 5854175.5
$ws.Range("I131").Value = 6789.4116
$ws.Range("J131").Value = 55556956
$ws.Range("K131").Value = 20368.2348
$ws.Range("L131").Value = 166670868
$ws.Range("M131").Value = -15328.2348
$ws.Range("N131").Value = -166680948

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 680.86957
$ws.Range("I107").Value = 558
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 558
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 1362
$ws.Range("N107").Value = -5340
$ws.Range("H113").Value = 25001940
$ws.Range("I113").Value = 50001220
$ws.Range("J113").Value = 2660
$ws.Range("K113").Value = 50001220
$ws.Range("L113").Value = 2660
$ws.Range("M113").Value = -49999050
$ws.Range("N113").Value = -7000
$ws.Range("H136").Value = 18433.334
$ws.Range("J136").Value = 18433.334
$ws.Range("L136").Value = 55300.00199999999
$ws.Range("N136").Value = -60400.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 35663
$ws.Range("J133").Value = 35663
$ws.Range("L133").Value = 35663
$ws.Range("N133").Value = -40723

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 78953
$ws.Range("J46").Value = 78953
$ws.Range("L46").Value = 78953
$ws.Range("N46").Value = -79415
$ws.Range("H107").Value = 1104.3077
$ws.Range("I107").Value = 1245.6
$ws.Range("K107").Value = 3736.8
$ws.Range("M107").Value = -1816.8
$ws.Range("H113").Value = 1056.8889
$ws.Range("I113").Value = 1280.2858
$ws.Range("J113").Value = 275
$ws.Range("K113").Value = 3840.8574
$ws.Range("L113").Value = 825
$ws.Range("M113").Value = -1670.8574
$ws.Range("N113").Value = -5165
$ws.Range("H134").Value = 78953
$ws.Range("J134").Value = 78953
$ws.Range("L134").Value = 236859
$ws.Range("N134").Value = -241929
$ws.Range("H136").Value = 9998.154
$ws.Range("I136").Value = 10802.167
$ws.Range("K136").Value = 32406.501
$ws.Range("M136").Value = -29856.501
$ws.Range("H137").Value = 40000
$ws.Range("J137").Value = 40000
$ws.Range("L137").Value = 40000
$ws.Range("N137").Value = -50200
